# Update the "想去人数" (F column) counts across the sheets to match
# the newly scraped totals (commit: "Update gh-pages to output generated
# at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 10
$ws.Range("F5").Value = 508
$ws.Range("F6").Value = 956
$ws.Range("F7").Value = 191
$ws.Range("F9").Value = 1024
$ws.Range("F10").Value = 819
$ws.Range("F11").Value = 239
$ws.Range("F14").Value = 824
$ws.Range("F15").Value = 280
$ws.Range("F18").Value = 1328
$ws.Range("F21").Value = 1182
$ws.Range("F22").Value = 2861
$ws.Range("F23").Value = 1409
$ws.Range("F25").Value = 191
$ws.Range("F28").Value = 1012
$ws.Range("F30").Value = 3073
$ws.Range("F31").Value = 595
$ws.Range("F33").Value = 1394

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 1

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 737

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 10
$ws.Range("F3").Value = 737
$ws.Range("F7").Value = 508
$ws.Range("F12").Value = 956
$ws.Range("F13").Value = 191
$ws.Range("F16").Value = 1024
$ws.Range("F17").Value = 819
$ws.Range("F18").Value = 239
$ws.Range("F26").Value = 824
$ws.Range("F27").Value = 280
$ws.Range("F30").Value = 1328
$ws.Range("F33").Value = 1182
$ws.Range("F34").Value = 2861
$ws.Range("F35").Value = 1409
$ws.Range("F37").Value = 191
$ws.Range("F42").Value = 1012
$ws.Range("F44").Value = 3073
$ws.Range("F45").Value = 595
$ws.Range("F47").Value = 1394
$ws.Range("F48").Value = 1
